$d = $word.ActiveDocument

# Locate the "Location (spatial scale)," phrase inside the "Key words/phrases" bullet
# and replace it with "Texas grasslands,". Touching the Font on just this sub-range
# (instead of the whole paragraph) causes Word to split the original single run into
# three runs: the unchanged lead-in text, the newly replaced phrase, and the unchanged
# trailing text -- matching how Word behaves when a mid-run selection is retyped.
$target = $d.Content.Duplicate
$target.Find.Execute("Location (spatial scale),", $true, $false, $false, $false, $false, $true, 1, $false, "Texas grasslands,", 2)

# Re-apply the (unchanged) run formatting explicitly on the replaced sub-range so Word
# materializes it as its own run instead of silently merging back into its neighbors.
$target.Font.Name = "Times New Roman"
$target.Font.NameBi = "Times New Roman"
$target.Font.Size = 12
